# Power_Storage.xlsx edit: "Increase MaxInvest Storage Adapt Szenarios Existing Units"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MaxInvest (column S) increased from 8 to 15 for the 5 BESS rows (7-11)
$ws.Range("S7:S11").Value = 15

# ExisUnits (column E) for row 10 (Node_6 / BESS7) increased from 30 to 33
$ws.Range("E10").Value = 33

# Update the on-screen selection to match what was selected after the edit
$ws.Range("S8:S11").Select() | Out-Null
